$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '95.533.02'
$ws.Range('E2').Value = '  +2.95%  '
$ws.Range('D3').Value = '3.594.67'
$ws.Range('E3').Value = '  +5.33%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '653.63'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.68%  '
$ws.Range('E7').Value = '  +7.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.407'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.17%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('E10').Value = '  +4.55%  '
$ws.Range('D11').Value = '3.591.68'
$ws.Range('E11').Value = '  +5.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.12'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.19%  '
$ws.Range('E13').Value = '  +1.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.32'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.46%  '
$ws.Range('D15').Value = '4.285.43'
$ws.Range('E15').Value = '  +5.99%  '
$ws.Range('D16').Value = '95.529.55'
$ws.Range('E16').Value = '  +3.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000256'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.32%  '
$ws.Range('D18').Value = '3.599.88'
$ws.Range('E18').Value = '  +5.62%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.94'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.58'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +9.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.53'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.92%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.485'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +12.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '511.12'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.96%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000196'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.64'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.70%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '96.87'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.75'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.39%  '
$ws.Range('D29').Value = '3.799.65'
$ws.Range('E29').Value = '  +5.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.20'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +18.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.33'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.03%  '
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.139'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.43%  '
$ws.Range('E34').Value = '  +1.09%  '
$ws.Range('E35').Value = '  +2.87%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '31.85'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.72%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.560'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.28'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +11.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '568.14'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.49'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.36%  '
$ws.Range('E41').Value = '  +1.53%  '
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('E43').Value = '  +0.50%  '
$ws.Range('E44').Value = '  +1.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.73'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.76'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.46%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0418'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.46%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.25'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.21%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '33.63'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +31.16%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '54.32'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.71%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.45'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.50%  '
